# Update "Last Updated" timestamp on the Metadata sheet.
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "30 Oct 2025, 01:12 AM"

# Insert a new row at the top of the "Top Gainers" ranking (row 36) for a
# newly-qualifying stock (SKMEGGPROD), pushing the existing rows 36-75 down
# to 37-76 and dropping the stock that falls off the bottom of the list
# (old row 76, CGPOWER).
$gainers = $wb.Worksheets.Item("Top Gainers")
$gainers.Rows.Item(36).Insert()

$gainers.Cells.Item(36, 1).Value = "🚀"
$gainers.Cells.Item(36, 2).Value = "SKMEGGPROD"
$gainers.Cells.Item(36, 3).Value = 4.9959
$gainers.Cells.Item(36, 4).Value = 6.6906
$gainers.Cells.Item(36, 5).Value = 23.7638

# The insert pushed the last row (old CGPOWER row) past the bottom of the
# table to row 77; remove it so the sheet ends at row 76 again.
$gainers.Rows.Item(77).Delete()
